# Update "想去人数" (F) and "最低票价" (G) figures on the 展览 and 全部类型 sheets.
# The same set of row edits is applied identically to both sheets, since
# their data is mirrored.

$wb = $excel.ActiveWorkbook

# Map of F-column updates per row number (applies to both sheets)
$fUpdates = @{
    2  = 11661
    3  = 11242
    4  = 604
    6  = 1018
    8  = 70
    9  = 43
    11 = 10715
    12 = 4142
    14 = 3
    15 = 10
    18 = 48
    19 = 122
    20 = 441
    21 = 11127
    22 = 10904
    24 = 27
}

# G-column updates per row number (applies to both sheets)
$gUpdates = @{
    20 = 49.9
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $fUpdates.Keys) {
        $ws.Range("F$row").Value = $fUpdates[$row]
    }

    foreach ($row in $gUpdates.Keys) {
        $ws.Range("G$row").Value = $gUpdates[$row]
    }
}
